$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.521.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.95%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.443.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.78%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'579.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.97%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'149.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +9.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.444.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.84%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.96%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.90%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.30%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.035.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.443.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.71%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'61.596.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +8.31%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.95%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.36%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'388.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +2.66%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.586.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'72.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.51%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.181"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.56%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'7.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -13.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'8.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.08%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.64%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D36").Value = "'24.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.98%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'7.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.87%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.53%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'166.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0793"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.41%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'26.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +9.85%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.87%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Filecoin"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'4.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'42.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.59%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.611.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.26%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'23.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.53%  "
$ws.Range("E51").Style = "Normal"
